# Apply updated numeric values to Sheet1 (Call_CCF data) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 152
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 128
$ws.Range("E2").Value = 98
$ws.Range("F2").Value = 155
$ws.Range("G2").Value = 226
$ws.Range("H2").Value = 190
$ws.Range("I2").Value = 63
$ws.Range("J2").Value = 91

# Row 3
$ws.Range("B3").Value = 26
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = 21
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 36
$ws.Range("G3").Value = 35
$ws.Range("H3").Value = 36
$ws.Range("I3").Value = 19
$ws.Range("J3").Value = 19

# Row 4
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 2

# Row 5
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 1

# Row 6
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 18
$ws.Range("I6").Value = 9
$ws.Range("J6").Value = 6

# Row 8
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 11
$ws.Range("F8").Value = 19
$ws.Range("G8").Value = 13
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 4

# Row 9
$ws.Range("F9").Value = 9
$ws.Range("H9").Value = 3
$ws.Range("J9").Value = 3

# Row 10
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 4

# Row 11
$ws.Range("B11").Value = 26
$ws.Range("C11").Value = 32
$ws.Range("D11").Value = 23
$ws.Range("E11").Value = 43
$ws.Range("F11").Value = 52
$ws.Range("G11").Value = 53
$ws.Range("H11").Value = 43
$ws.Range("I11").Value = 7
$ws.Range("J11").Value = 19

# Row 12
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 2

# Row 13
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 6
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1

# Row 14
$ws.Range("B14").Value = 6
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 9
$ws.Range("H14").Value = 9
$ws.Range("J14").Value = 2

# Row 15
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 6
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 5

# Row 16
$ws.Range("B16").Value = 7
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 17
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 15
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 9

# Row 17
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 18
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 7
$ws.Range("J17").Value = 2

# Row 18
$ws.Range("B18").Value = 16
$ws.Range("C18").Value = 11
$ws.Range("D18").Value = 18
$ws.Range("E18").Value = 6
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = 33
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 10

# Row 19
$ws.Range("B19").Value = 4
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 1
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 6
$ws.Range("J19").Value = 2

# Row 20
$ws.Range("B20").Value = 7
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 13
$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 4

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 13
$ws.Range("H22").Value = 8
$ws.Range("J22").Value = 4

# Row 23
$ws.Range("B23").Value = 5
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 6
$ws.Range("I23").Value = 2

# Row 24
$ws.Range("B24").Value = 60
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 47
$ws.Range("E24").Value = 26
$ws.Range("F24").Value = 41
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 55
$ws.Range("I24").Value = 23
$ws.Range("J24").Value = 28

# Row 25
$ws.Range("B25").Value = 4
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 6
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 5
$ws.Range("I25").Value = 1

# Row 26
$ws.Range("B26").Value = 8
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 1
$ws.Range("G26").Value = 5
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = 2

# Row 27
$ws.Range("B27").Value = 7
$ws.Range("C27").Value = 9
$ws.Range("D27").Value = 6
$ws.Range("E27").Value = 3
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = 12
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 4

# Row 28
$ws.Range("B28").Value = 14
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 5
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 20
$ws.Range("H28").Value = 13
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 12

# Row 29
$ws.Range("B29").Value = 6
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 10
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = 13
$ws.Range("H29").Value = 7
$ws.Range("J29").Value = 4

# Row 30
$ws.Range("B30").Value = 18
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 10
$ws.Range("E30").Value = 2
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 15
$ws.Range("H30").Value = 5
$ws.Range("I30").Value = 9
$ws.Range("J30").Value = 5

# Row 31
$ws.Range("B31").Value = 3
$ws.Range("D31").Value = 7
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 13
$ws.Range("H31").Value = 13
$ws.Range("J31").Value = 1

# Row 32
$ws.Range("B32").Value = 16
$ws.Range("C32").Value = 9
$ws.Range("D32").Value = 13
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = 12
$ws.Range("G32").Value = 17
$ws.Range("H32").Value = 13
$ws.Range("I32").Value = 3
$ws.Range("J32").Value = 12

# Row 33
$ws.Range("B33").Value = 5
$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 9
$ws.Range("E33").Value = 12
$ws.Range("F33").Value = 2
$ws.Range("G33").Value = 11
$ws.Range("H33").Value = 3
$ws.Range("I33").Value = 2
$ws.Range("J33").Value = 1

# Row 34
$ws.Range("B34").Value = 5
$ws.Range("D34").Value = 2
$ws.Range("F34").Value = 7
$ws.Range("G34").Value = 4
$ws.Range("H34").Value = 5
$ws.Range("I34").Value = 1
$ws.Range("J34").Value = 7

# Row 35
$ws.Range("B35").Value = 6
$ws.Range("C35").Value = 3
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = 1
$ws.Range("F35").Value = 3
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 5
$ws.Range("J35").Value = 4

# Row 38
$ws.Range("B38").Value = 8
$ws.Range("D38").Value = 6
$ws.Range("E38").Value = 2
$ws.Range("F38").Value = 2
$ws.Range("G38").Value = 9
$ws.Range("H38").Value = 10
$ws.Range("I38").Value = 7
$ws.Range("J38").Value = 3

# Row 39
$ws.Range("B39").Value = 6
$ws.Range("C39").Value = 2
$ws.Range("D39").Value = 4
$ws.Range("G39").Value = 4
$ws.Range("H39").Value = 7
$ws.Range("I39").Value = 7
$ws.Range("J39").Value = 3

# Row 42
$ws.Range("B42").Value = 2
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 2
$ws.Range("E42").Value = 2
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 5
$ws.Range("H42").Value = 3
